$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 6 (pushes the old rows 6-10 down to 7-11)
$ws.Rows.Item(6).Insert()

# Content for the newly inserted row
$ws.Range("A6").Value = "被"
$ws.Range("B6").Value = "passive"
$ws.Range("C6").Value = "虛詞"

# Make rows 5 and 6 a bit taller (matches the thicker-bottom-border look)
$ws.Rows.Item(5).RowHeight = 16.5
$ws.Rows.Item(6).RowHeight = 16.5

# Draw a red medium "box" across A6:C6 - left edge only on A6, right edge
# only on C6, and a top+bottom edge shared by every cell in between.
$rowRange = $ws.Range("A6:C6")
$rowRange.BorderAround([System.Reflection.Missing]::Value, -4138, [System.Reflection.Missing]::Value, 255)

# Red font for the new row's text
$rowRange.Font.Color = 255

# Match the authored selection state
$ws.Range("A6:C6").Select() | Out-Null
